# refactor(translations): create new translator
#
# The two translation worksheets are renamed from the old
# "src|assets\translations\..." scheme to a new per-app naming scheme
# (app1/app2), and the second sheet's "key" column entries are
# re-namespaced from "app.*" to "app2.*" now that a dedicated app2
# translator exists. The second sheet also becomes the active/selected
# sheet, matching the new authoring workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "src|assets\translations\transla" -> "app1|app1.json" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "app1|app1.json"

# --- Sheet 2: "src|assets\translations\app\app" -> "app2|app2.json" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "app2|app2.json"

# Re-key the translation keys used on the second ("app2") sheet.
$ws2.Range("A2").Value = "app2.translation"
$ws2.Range("A3").Value = "app2.title"

# Make the second sheet the active / visible tab, with B9 selected.
$ws2.Activate()
[void]$ws2.Range("B9").Select()
